# Update the "Playcount" figures for a few artists (C2, C6, C10, C11).
# The source values are stored as text (not numbers), so the cell's
# number format is forced to Text ("@") before assigning the new
# string value - this prevents Excel from reinterpreting the
# digit-only string as a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2"  = "1033"
    "C6"  = "100"
    "C10" = "71"
    "C11" = "69"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
